# Update the ammonium.N (column G) values on the "slurry" sheet with the
# newly-analysed TAN figures, applying the same 2-decimal number format
# used elsewhere in the table (e.g. TS/VS columns D:E).
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("slurry")
$ws2 = $wb.Worksheets.Item("soil")

$ws1.Range("G2").Value = 2.008686210640608
$ws1.Range("G3").Value = 1.9821605550049555
$ws1.Range("G4").Value = 1.8901890189018902
$ws1.Range("G5").Value = 1.7690875232774672
$ws1.Range("G6").Value = 1.875
$ws1.Range("G7").Value = 1.9762845849802373
$ws1.Range("G2:G7").NumberFormat = "0.00"

# Bring the "slurry" sheet to the front (it was "soil" before) and move the
# selection to I14. "soil" sheet's own tab-selected state reverts to
# not-selected automatically since it is no longer the active sheet.
$null = $ws1.Activate()
$null = $ws1.Range("I14").Select()
